$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 1.02),
    @(2, 3, 1.061049720912395),
    @(2, 4, 1.066633807911769),
    @(2, 5, 1.056691385512943),
    @(2, 6, 1.074611611373777),
    @(2, 9, 1.04698292226123),
    @(2, 10, 1.066027408165287),
    @(2, 11, 1.069343830705615),
    @(2, 12, 1.059428401738694),
    @(2, 13, 1.077300365105973),
    @(2, 14, 1.067541289359935),
    @(3, 2, 1.02),
    @(3, 3, 1.062556202126101),
    @(3, 4, 1.068007345698905),
    @(3, 5, 1.057993591883989),
    @(3, 6, 1.07606399533798),
    @(3, 9, 1.047372343799747),
    @(3, 10, 1.067185121795014),
    @(3, 11, 1.070531607977588),
    @(3, 12, 1.06054305454678),
    @(3, 13, 1.078568348806519),
    @(3, 14, 1.068700647075808),
    @(4, 2, 1.02),
    @(4, 3, 1.063529455790313),
    @(4, 4, 1.068894872495675),
    @(4, 5, 1.058835034177963),
    @(4, 6, 1.077002610832716),
    @(4, 9, 1.047621890076945),
    @(4, 10, 1.067932279627217),
    @(4, 11, 1.071298400105542),
    @(4, 12, 1.061262590693808),
    @(4, 13, 1.079387122754626),
    @(4, 14, 1.069448865957832),
    @(5, 2, 1.02),
    @(5, 3, 1.063938250358308),
    @(5, 4, 1.069267697061578),
    @(5, 5, 1.059188501559345),
    @(5, 6, 1.077396930139711),
    @(5, 9, 1.047726218491433),
    @(5, 10, 1.068245921063046),
    @(5, 11, 1.071620339195516),
    @(5, 12, 1.061564677479499),
    @(5, 13, 1.079730935090483),
    @(5, 14, 1.069762952800476),
    @(6, 2, 1.02),
    @(6, 3, 1.064006867797172),
    @(6, 4, 1.069330279018318),
    @(6, 5, 1.059247834329063),
    @(6, 6, 1.077463122166894),
    @(6, 9, 1.047743701694726),
    @(6, 10, 1.068298555809172),
    @(6, 11, 1.071674369700194),
    @(6, 12, 1.061615375511994),
    @(6, 13, 1.079788639394416),
    @(6, 14, 1.069815662293979),
    @(7, 2, 1.02),
    @(7, 3, 1.063534919533614),
    @(7, 4, 1.068899855334395),
    @(7, 5, 1.058839758299778),
    @(7, 6, 1.077007880817835),
    @(7, 9, 1.047623286397908),
    @(7, 10, 1.06793647233464),
    @(7, 11, 1.071302703516056),
    @(7, 12, 1.061266628781858),
    @(7, 13, 1.079391718357288),
    @(7, 14, 1.069453064619381),
    @(8, 2, 1.02),
    @(8, 3, 1.061559166776242),
    @(8, 4, 1.067098262708112),
    @(8, 5, 1.057131718423172),
    @(8, 6, 1.075102698306535),
    @(8, 9, 1.047115034723435),
    @(8, 10, 1.06641907217128),
    @(8, 11, 1.069745617383267),
    @(8, 12, 1.059805462930798),
    @(8, 13, 1.077729240904035),
    @(8, 14, 1.06793350957374),
    @(9, 2, 1.02),
    @(9, 3, 1.05806550271174),
    @(9, 4, 1.063913823573324),
    @(9, 5, 1.054112699313346),
    @(9, 6, 1.0717362424809),
    @(9, 9, 1.046200676852647),
    @(9, 10, 1.063729954301296),
    @(9, 11, 1.066987953181933),
    @(9, 12, 1.057217301570183),
    @(9, 13, 1.074786487653932),
    @(9, 14, 1.065240572847935),
    @(10, 2, 1.02),
    @(10, 3, 1.055727759974409),
    @(10, 4, 1.061783885965812),
    @(10, 5, 1.05209345045221),
    @(10, 6, 1.069485291888571),
    @(10, 9, 1.045578353757689),
    @(10, 10, 1.061926591817482),
    @(10, 11, 1.065139822381255),
    @(10, 12, 1.055482517547885),
    @(10, 13, 1.072815368109503),
    @(10, 14, 1.063434649382545),
    @(11, 2, 1.02),
    @(11, 3, 1.054713337615693),
    @(11, 4, 1.060859857379696),
    @(11, 5, 1.051217454248732),
    @(11, 6, 1.068508937603825),
    @(11, 9, 1.045305824305496),
    @(11, 10, 1.061143118167658),
    @(11, 11, 1.064337183485568),
    @(11, 12, 1.054729045986332),
    @(11, 13, 1.071959566030656),
    @(11, 14, 1.062650063110274),
    @(12, 2, 1.02),
    @(12, 3, 1.054336201579225),
    @(12, 4, 1.060516361287374),
    @(12, 5, 1.05089181558003),
    @(12, 6, 1.068146016114691),
    @(12, 9, 1.045204132222783),
    @(12, 10, 1.060851702731177),
    @(12, 11, 1.064038682198271),
    @(12, 12, 1.054448821218485),
    @(12, 13, 1.071641331048326),
    @(12, 14, 1.062358233830446),
    @(13, 2, 1.02),
    @(13, 3, 1.054417113820939),
    @(13, 4, 1.060590054738858),
    @(13, 5, 1.050961677818475),
    @(13, 6, 1.068223875901765),
    @(13, 9, 1.045225966495537),
    @(13, 10, 1.0609142304739),
    @(13, 11, 1.064102728388458),
    @(13, 12, 1.054508946411355),
    @(13, 13, 1.071709609578127),
    @(13, 14, 1.062420850369736),
    @(14, 2, 1.02),
    @(14, 3, 1.054682170308745),
    @(14, 4, 1.060831469452509),
    @(14, 5, 1.051190542091149),
    @(14, 6, 1.068478943751255),
    @(14, 9, 1.045297427857939),
    @(14, 10, 1.0611190378452),
    @(14, 11, 1.064312516772006),
    @(14, 12, 1.054705889735892),
    @(14, 13, 1.07193326785054),
    @(14, 14, 1.062625948590994),
    @(15, 2, 1.02),
    @(15, 3, 1.054845435769244),
    @(15, 4, 1.060980176922607),
    @(15, 5, 1.051331518963821),
    @(15, 6, 1.068636064775087),
    @(15, 9, 1.045341396188665),
    @(15, 10, 1.061245173420304),
    @(15, 11, 1.064441725681536),
    @(15, 12, 1.054827186193865),
    @(15, 13, 1.072071024233131),
    @(15, 14, 1.062752263293083),
    @(16, 2, 1.02),
    @(16, 3, 1.055795037458577),
    @(16, 4, 1.061845173150958),
    @(16, 5, 1.052151552109604),
    @(16, 6, 1.069550053279026),
    @(16, 9, 1.045596375937192),
    @(16, 10, 1.061978532890542),
    @(16, 11, 1.065193040014843),
    @(16, 12, 1.055532473987546),
    @(16, 13, 1.072872115915524),
    @(16, 14, 1.063486664217887),
    @(17, 2, 1.02),
    @(17, 3, 1.056390111558),
    @(17, 4, 1.062387288159072),
    @(17, 5, 1.052665491019124),
    @(17, 6, 1.070122919318725),
    @(17, 9, 1.045755496884469),
    @(17, 10, 1.062437847379182),
    @(17, 11, 1.065663675916361),
    @(17, 12, 1.055974262378615),
    @(17, 13, 1.073374000182777),
    @(17, 14, 1.063946630985785),
    @(18, 2, 1.02),
    @(18, 3, 1.056737000014648),
    @(18, 4, 1.0627033261887),
    @(18, 5, 1.0529651043508),
    @(18, 6, 1.070456901043731),
    @(18, 9, 1.045848014444083),
    @(18, 10, 1.06270550680103),
    @(18, 11, 1.065937960130102),
    @(18, 12, 1.056231728732958),
    @(18, 13, 1.073666519913499),
    @(18, 14, 1.064214670514713),
    @(19, 2, 1.02),
    @(19, 3, 1.056855244978174),
    @(19, 4, 1.062811058640217),
    @(19, 5, 1.053067238075744),
    @(19, 6, 1.070570753053933),
    @(19, 9, 1.045879510590461),
    @(19, 10, 1.062796729390768),
    @(19, 11, 1.066031445200262),
    @(19, 12, 1.056319480759938),
    @(19, 13, 1.073766224338084),
    @(19, 14, 1.064306022650992),
    @(20, 2, 1.02),
    @(20, 3, 1.056326287340671),
    @(20, 4, 1.062329141831161),
    @(20, 5, 1.052610366687657),
    @(20, 6, 1.07006147296048),
    @(20, 9, 1.045738455246389),
    @(20, 10, 1.062388593288589),
    @(20, 11, 1.065613204952427),
    @(20, 12, 1.055926885606899),
    @(20, 13, 1.073320175649283),
    @(20, 14, 1.063897306948733),
    @(21, 2, 1.02),
    @(21, 3, 1.054604127085237),
    @(21, 4, 1.06076038637121),
    @(21, 5, 1.051123154354429),
    @(21, 6, 1.068403839869502),
    @(21, 9, 1.045276397058508),
    @(21, 10, 1.06105873823399),
    @(21, 11, 1.064250749466866),
    @(21, 12, 1.054647904601925),
    @(21, 13, 1.071867415839243),
    @(21, 14, 1.062565563347421),
    @(22, 2, 1.02),
    @(22, 3, 1.053519397184751),
    @(22, 4, 1.059772477370504),
    @(22, 5, 1.050186607921747),
    @(22, 6, 1.067360111676771),
    @(22, 9, 1.04498320535873),
    @(22, 10, 1.060220296851691),
    @(22, 11, 1.063392001664018),
    @(22, 12, 1.053841719326158),
    @(22, 13, 1.07095196795015),
    @(22, 14, 1.061725931282189),
    @(23, 2, 1.02),
    @(23, 3, 1.054094619964086),
    @(23, 4, 1.060296338005646),
    @(23, 5, 1.050683231281858),
    @(23, 6, 1.06791355740035),
    @(23, 9, 1.045138886573602),
    @(23, 10, 1.060664991667118),
    @(23, 11, 1.063847443227165),
    @(23, 12, 1.054269288878772),
    @(23, 13, 1.07143746002803),
    @(23, 14, 1.06217125761526),
    @(24, 2, 1.02),
    @(24, 3, 1.056355127400746),
    @(24, 4, 1.062355416174731),
    @(24, 5, 1.052635275490384),
    @(24, 6, 1.070089238423664),
    @(24, 9, 1.045746156540722),
    @(24, 10, 1.062410849867646),
    @(24, 11, 1.065636011318656),
    @(24, 12, 1.05594829381487),
    @(24, 13, 1.073344497321843),
    @(24, 14, 1.063919595134685),
    @(25, 2, 1.02),
    @(25, 3, 1.058970185319039),
    @(25, 4, 1.06473827878961),
    @(25, 5, 1.054894319422229),
    @(25, 6, 1.072607693514297),
    @(25, 9, 1.04643929723024),
    @(25, 10, 1.064427000963105),
    @(25, 11, 1.067702556768185),
    @(25, 12, 1.0578880255063),
    @(25, 13, 1.075548868846291),
    @(25, 14, 1.065938609395956)
)

foreach ($row in $data) {
    $ws.Cells.Item($row[0], $row[1]).Value = $row[2]
}
